$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Greedy results block (P1:Q18 -> P1:Q19) ---
# Row 18 currently holds the "21 cities" data point; a new "20 cities" data
# point is inserted before it, so the old row 18 values move to row 19 and
# the new data point takes row 18.
$ws.Range("P19").Value = 21
$ws.Range("Q19").Value = 0.48313867179154701
$ws.Range("S19").Value = 0.00098596230236580595

$ws.Range("P18").Value = 20
$ws.Range("Q18").Value = 0.55610751928873203
$ws.Range("S18").Value = 0.00036107909982092598

# --- MST results block (P21:Q38 -> P21:Q39) ---
# Same pattern: row 38 currently holds the "21 cities" data point; the new
# "20 cities" data point is inserted before it.
$ws.Range("P39").Value = 21
$ws.Range("Q39").Value = 1.2346367360332
$ws.Range("S39").Value = 0.00021869589836569499

$ws.Range("P38").Value = 20
$ws.Range("Q38").Value = 1.34038524510091
$ws.Range("S38").Value = 0.000190300005488097

# --- Update defined names to cover the expanded ranges ---
$wb.Names.Item("result_greedy").RefersTo = "=Sheet1!`$P`$1:`$Q`$19"
$wb.Names.Item("result_mst").RefersTo = "=Sheet1!`$P`$21:`$Q`$39"
